$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date string in column A (A2:A6) from 2025-11-29 to 2025-12-01
# Force text format so Excel does not auto-convert the string into a date serial value,
# then restore the original (default) cell style so no stray style index is left on the cells.
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("A2:A6").Value = "2025-12-01"
$ws.Range("A2:A6").Style = "Normal"

# Update the N column values (N2:N6) from 85.87127175646313 to 85.87246918135976
$ws.Range("N2:N6").Value = 85.87246918135976
